$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.782.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.086.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.652"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.31"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.32"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.369"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0768"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.97%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.95"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.883"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.389.26"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.109.64"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.746.71"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.08%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.49"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.95"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.27%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.10"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "21.12"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.28%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.124"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.32"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +10.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.70"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0611"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0838"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.28"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.03%  "

$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.94"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.92%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.16"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0221"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.93"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.15%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.63"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.337.27"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.44"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.02"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.14%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.272.96"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.06%  "

Write-Output "Updated cryptos list"